$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @('ECs', 'C3', 'Itgam', 'FAPs', 3, 1, 0.5236730000000001, 1.571019, 0.001411687926676084, 0.001412631582657676, 3, 1, 0.9636300000000001, 2.89089, 0.00463428088190967, 0.004636691828827675, 0.5046270129900001, 4.54164311691, 0.000006542158369817673, 0.000006549937316452754),
  @('ECs', 'C3', 'Itgam', 'M1', 3, 1, 0.5236730000000001, 1.571019, 0.001411687926676084, 0.001412631582657676, 3, 1, 71.62512099999999, 214.875363, 0.3444588990740914, 0.344638101013349, 37.508141989433, 337.573277904897, 0.0004862684690590305, 0.0004868466660786233),
  @('ECs', 'C3', 'Itgam', 'M2', 3, 1, 0.5236730000000001, 1.571019, 0.001411687926676084, 0.001412631582657676, 3, 1, 135.0220616666666, 405.066185, 0.6493469059886695, 0.6496847234325412, 70.70740809916833, 636.366672892515, 0.0009166751874086746, 0.0009177651590910254),
  @('ECs', 'C3', 'Itgam', 'sCs', 3, 1, 0.5236730000000001, 1.571019, 0.001411687926676084, 0.001412631582657676, 1, 0.5, 0.324361, 0.648722, 0.001559914055329434, 0.001040483725282092, 0.169859097953, 1.019154587718, 0.00000220211183856089, 0.000001469820171574796),
  @('FAPs', 'C3', 'Itgam', 'FAPs', 3, 1, 290.8976950000001, 872.6930850000001, 0.7841854820267643, 0.78470967813754, 3, 1, 0.9636300000000001, 2.89089, 0.00463428088190967, 0.004636691828827675, 280.3177458328501, 2522.85971249565, 0.003634135787227753, 0.003638456952622326),
  @('FAPs', 'C3', 'Itgam', 'M1', 3, 1, 290.8976950000001, 872.6930850000001, 0.7841854820267643, 0.78470967813754, 3, 1, 71.62512099999999, 214.875363, 0.3444588990740914, 0.344638101013349, 20835.5826029961, 187520.2434269649, 0.270119667808825, 0.2704408533201181),
  @('FAPs', 'C3', 'Itgam', 'M2', 3, 1, 290.8976950000001, 872.6930850000001, 0.7841854820267643, 0.78470967813754, 3, 1, 135.0220616666666, 405.066185, 0.6493469059886695, 0.6496847234325412, 39277.60651298119, 353498.4586168308, 0.5092084164753128, 0.5098138902156262),
  @('FAPs', 'C3', 'Itgam', 'sCs', 3, 1, 290.8976950000001, 872.6930850000001, 0.7841854820267643, 0.78470967813754, 1, 0.5, 0.324361, 0.648722, 0.001559914055329434, 0.001040483725282092, 94.35586724789502, 566.1352034873701, 0.001223261955398837, 0.000816477649173459),
  @('M1', 'C3', 'Itgam', 'FAPs', 3, 1, 38.68424166666667, 116.052725, 0.1042827812651277, 0.1043524900644015, 3, 1, 0.9636300000000001, 2.89089, 0.00463428088190967, 0.004636691828827675, 37.27729579725001, 335.49566217525, 0.0004832756995293492, 0.0004838503379994317),
  @('M1', 'C3', 'Itgam', 'M1', 3, 1, 38.68424166666667, 116.052725, 0.1042827812651277, 0.1043524900644015, 3, 1, 71.62512099999999, 214.875363, 0.3444588990740914, 0.344638101013349, 2770.763490168241, 24936.87141151417, 0.03592113202697018, 0.03596384401180971),
  @('M1', 'C3', 'Itgam', 'M2', 3, 1, 38.68424166666667, 116.052725, 0.1042827812651277, 0.1043524900644015, 3, 1, 135.0220616666666, 405.066185, 0.6493469059886695, 0.6496847234325412, 5223.226063844902, 47009.03457460412, 0.06771570136240386, 0.06779621864698772),
  @('M1', 'C3', 'Itgam', 'sCs', 3, 1, 38.68424166666667, 116.052725, 0.1042827812651277, 0.1043524900644015, 1, 0.5, 0.324361, 0.648722, 0.001559914055329434, 0.001040483725282092, 12.54765931124167, 75.28595586745, 0.0001626721762243177, 0.000108577067604671),
  @('M2', 'C3', 'Itgam', 'FAPs', 3, 1, 40.106198, 120.318594, 0.1081160103757125, 0.1081882815327926, 3, 1, 0.9636300000000001, 2.89089, 0.00463428088190967, 0.004636691828827675, 38.64753557874, 347.82782020866, 0.0005010399599125118, 0.0005016357209580077),
  @('M2', 'C3', 'Itgam', 'M1', 3, 1, 40.106198, 120.318594, 0.1081160103757125, 0.1081882815327926, 3, 1, 71.62512099999999, 214.875363, 0.3444588990740914, 0.344638101013349, 2872.611284599958, 25853.50156139962, 0.03724152190630096, 0.03728580389935923),
  @('M2', 'C3', 'Itgam', 'M2', 3, 1, 40.106198, 120.318594, 0.1081160103757125, 0.1081882815327926, 3, 1, 135.0220616666666, 405.066185, 0.6493469059886695, 0.6496847234325412, 5415.221539571542, 48736.99385614388, 0.07020479682530778, 0.07028827376627429),
  @('M2', 'C3', 'Itgam', 'sCs', 3, 1, 40.106198, 120.318594, 0.1081160103757125, 0.1081882815327926, 1, 0.5, 0.324361, 0.648722, 0.001559914055329434, 0.001040483725282092, 13.008886489478, 78.053318936868, 0.0001686516841912168, 0.0001125681462011078),
  @('sCs', 'C3', 'Itgam', 'FAPs', 2, 1, 0.7434085, 1.486817, 0.002004038405719556, 0.001336918682608127, 3, 1, 0.9636300000000001, 2.89089, 0.00463428088190967, 0.004636691828827675, 0.7163707328550001, 4.29822439713, 0.000009287276870238873, 0.000006198879931456165),
  @('sCs', 'C3', 'Itgam', 'M1', 2, 1, 0.7434085, 1.486817, 0.002004038405719556, 0.001336918682608127, 3, 1, 71.62512099999999, 214.875363, 0.3444588990740914, 0.344638101013349, 53.2467237649285, 319.480342589571, 0.0006903088629363557, 0.0004607531159833333),
  @('sCs', 'C3', 'Itgam', 'M2', 2, 1, 0.7434085, 1.486817, 0.002004038405719556, 0.001336918682608127, 3, 1, 135.0220616666666, 405.066185, 0.6493469059886695, 0.6496847234325412, 100.3765483305242, 602.259289983145, 0.00130131613823646, 0.0008685756445620587),
  @('sCs', 'C3', 'Itgam', 'sCs', 2, 1, 0.7434085, 1.486817, 0.002004038405719556, 0.001336918682608127, 1, 0.5, 0.324361, 0.648722, 0.001559914055329434, 0.001040483725282092, 0.2411327244685, 0.964530897874, 0.000003126127676501926, 0.000001391042131279331)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}
